$wb = $excel.ActiveWorkbook

# --- Sheet: semantic_aspect_model_schema (header row) ---
$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")

$headerCols = @("O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI")
foreach ($col in $headerCols) {
    $cell = $wsSchema.Range($col + "1")
    $oldText = [string]$cell.Value()
    $cell.Value = $oldText.Replace("_", "__")
}

# Recompute the autofit column widths now that header text changed
$wsSchema.Range("O1:AI1").Columns.AutoFit() | Out-Null

# --- Sheet: description ---
$wsDesc = $wb.Worksheets.Item("description")

$wsDesc.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

$wsDesc.Range("B5").Value = "Digital Twin Field Name: id"
$wsDesc.Range("B6").Value = "Digital Twin Field Name: manufacturerPartId"
$wsDesc.Range("B7").Value = "Digital Twin Field Name: digitalTwinType"

$rowRange = 19..39
foreach ($r in $rowRange) {
    $cell = $wsDesc.Range("A" + $r)
    $oldText = [string]$cell.Value()
    $cell.Value = $oldText.Replace("_", "__")
}
